$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Corrected data entry errors ---
$ws.Cells.Item(85,4).Value = "B"
$ws.Cells.Item(207,3).Value = 37
$ws.Cells.Item(208,3).Value = 25
$ws.Cells.Item(209,3).Value = 20
$ws.Cells.Item(211,3).Value = 15

# --- New draft pick rows (261-284) ---
$rows = @(
    @{row=261; team="deano"; player="Jon Gray"; salary=12; pos="P"; date=43117},
    @{row=262; team="deano"; player="Starlin Castro"; salary=9; pos="2B"; date=43117},
    @{row=263; team="hermanos"; player="Max Kepler"; salary=14; pos="OF"; date=43117},
    @{row=264; team="hermanos"; player="Salvador Perez"; salary=20; pos="C"; date=43117},
    @{row=265; team="jobu"; player="Dylan Bundy"; salary=6; pos="P"; date=43117},
    @{row=266; team="jobu"; player="Joc Pederson"; salary=17; pos="OF"; date=43117},
    @{row=267; team="pasadena"; player="Yu Darvish"; salary=15; pos="P"; date=43117},
    @{row=268; team="rippe"; player="Daniel Murphy"; salary=19; pos="MI"; date=43117},
    @{row=269; team="balco"; player="Justus Sheffield"; salary=4; pos="P"; date=43118},
    @{row=270; team="balco"; player="Kevin Gausman"; salary=21; pos="P"; date=43118},
    @{row=271; team="chicago"; player="Shane Bieber"; salary=13; pos="P"; date=43118},
    @{row=272; team="deener"; player="Kyle Schwarber"; salary=18; pos="OF"; date=43118},
    @{row=273; team="deener"; player="Marcus Semien"; salary=9; pos="SS"; date=43118},
    @{row=274; team="hermanos"; player="Seranthony Dominguez"; salary=12; pos="P"; date=43118},
    @{row=275; team="hermanos"; player="Archie Bradley"; salary=12; pos="P"; date=43118},
    @{row=276; team="pkdodgers"; player="DJ LeMahieu"; salary=10; pos="MI"; date=43118},
    @{row=277; team="hermanos"; player="Justin Upton"; salary=26; pos="OF"; date=43121},
    @{row=278; team="hermanos"; player="Ender Inciarte"; salary=19; pos="OF"; date=43121},
    @{row=279; team="balco"; player="Yoan Moncada"; salary=15; pos="2B"; date=43121},
    @{row=280; team="chicago"; player="Adam Eaton"; salary=13; pos="OF"; date=43121},
    @{row=281; team="isotopes"; player="Carlos Santana"; salary=15; pos="CI"; date=43121},
    @{row=282; team="pasadena"; player="Byron Buxton"; salary=15; pos="OF"; date=43121},
    @{row=283; team="allrise"; player="David Robertson"; salary=8; pos="P"; date=43121},
    @{row=284; team="pkdodgers"; player="Rich Hill"; salary=15; pos="P"; date=43121}
)

$byPlayer = @{}
foreach ($r in $rows) {
    $byPlayer[$r.player] = $r
}

$playerOrder = @(
    "Jon Gray",
    "Starlin Castro",
    "Daniel Murphy",
    "Max Kepler",
    "Dylan Bundy",
    "Salvador Perez",
    "Joc Pederson",
    "Yu Darvish",
    "Justus Sheffield",
    "Seranthony Dominguez",
    "Kevin Gausman",
    "Kyle Schwarber",
    "Marcus Semien",
    "Shane Bieber",
    "DJ LeMahieu",
    "Archie Bradley",
    "Justin Upton",
    "Ender Inciarte",
    "Yoan Moncada",
    "Adam Eaton",
    "Carlos Santana",
    "Byron Buxton",
    "David Robertson",
    "Rich Hill"
)

# Copy formats from the last existing row (260) into the new row block so
# number formats (currency / date) match without creating new style entries.
$ws.Range("A260:E260").Copy()
$ws.Range("A261:E284").PasteSpecial(-4122)

# Seed the B-column (player) values in the exact order the names were first
# typed so the shared-string table is built up in that order.
foreach ($name in $playerOrder) {
    $rec = $byPlayer[$name]
    $ws.Cells.Item($rec.row, 2).Value = $rec.player
}

# Fill in the remaining columns per row.
foreach ($r in $rows) {
    $ws.Cells.Item($r.row, 1).Value = $r.team
    $ws.Cells.Item($r.row, 3).Value = $r.salary
    $ws.Cells.Item($r.row, 4).Value = $r.pos
    $ws.Cells.Item($r.row, 5).Value = $r.date
}

# --- Sort the first batch of new picks (rows 261-276) by draft date ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("E261:E276"))
$ws.Sort.SetRange($ws.Range("A261:E276"))
$ws.Sort.Apply()

# --- Restore selection/viewport close to the authored state ---
$ws.Range("A276:E284").Select()

Write-Host "done"